# Update the naive QoQ matched-to-ifoCAST error table (T45) with the
# "ifoCAST full series evaluation" figures for rows 2-11 (columns B:G).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @(0.1666076251416388, 1.464488876860186, 9.385379659755237, 3.063556700920556, 3.131009834800738, 22)
    3  = @(0.0708499967161877, 1.419337315890072, 9.02840922660916,  3.004731140486476, 3.078077156245957, 21)
    4  = @(-0.4288746198142434, 0.8403831407438848, 3.678417349840328, 1.917920058250689, 1.917916608230349, 20)
    5  = @(-0.01747259896522734, 0.5856938463863763, 0.8820486078156717, 0.9391744288552961, 0.9647429999147259, 19)
    6  = @(-0.01422969815280487, 0.6208678124689757, 0.9450450342577459, 0.9721342676080017, 1.000210739531721, 18)
    7  = @(-0.1197329514763359, 0.4398200131717107, 0.4468041506909826, 0.6684341034769117, 0.6778623988027417, 17)
    8  = @(-0.04262487491856422, 0.4573970767948405, 0.4014438679027904, 0.6335959816024643, 0.6528926306346485, 16)
    9  = @(0.04839913305653263, 0.3941245839882653, 0.3208696055700829, 0.566453533460674, 0.5841909510683425, 15)
    10 = @(-0.01748806567405372, 0.3260429352179946, 0.2494121353768383, 0.4994117893851108, 0.5179462580602099, 14)
    11 = @(-0.02168945670833724, 0.351174191232308, 0.2265631253025311, 0.4759864759659995, 0.4949078202060265, 13)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    for ($i = 0; $i -lt $values.Length; $i++) {
        $col = 2 + $i   # column B = 2
        $ws.Cells.Item($row, $col).Value = $values[$i]
    }
}
